# "first draft of postpartum event"
#
# The underlying edit collapses the two-value (B/C) pairs for several
# case-fatality-rate parameters down to a single value in column B,
# clearing column C entirely for those rows.
#
# Affected parameter rows on "parameter_values":
#   36  cfr_aph                    0.6 / 0.02   -> 0.02
#   37  cfr_eclampsia              0.5 / 0.184  -> 0.184
#   38  cfr_sepsis                 0.5 / 0.33   -> 0.33
#   39  cfr_uterine_rupture        0.8 / 0.345  -> 0.345
#   50  cfr_pp_pph                 0.5 / 0.1    -> 0.1
#   51  cfr_pp_eclampsia           0.5 / 0.184  -> 0.184
#   52  cfr_pp_sepsis              0.5 / 0.33   -> 0.33

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameter_values")
$ws.Activate()

# Row 36: cfr_aph -> keep only the (former column C) value, drop column C
$ws.Range("B36").Value = 0.02
$ws.Range("C36").ClearContents()

# Row 37: cfr_eclampsia
$ws.Range("B37").Value = 0.184
$ws.Range("C37").ClearContents()

# Row 38: cfr_sepsis
$ws.Range("B38").Value = 0.33
$ws.Range("C38").ClearContents()

# Row 39: cfr_uterine_rupture
$ws.Range("B39").Value = 0.34499999999999997
$ws.Range("C39").ClearContents()

# Row 50: cfr_pp_pph
$ws.Range("B50").Value = 0.1
$ws.Range("C50").ClearContents()

# Row 51: cfr_pp_eclampsia
$ws.Range("B51").Value = 0.184
$ws.Range("C51").ClearContents()

# Row 52: cfr_pp_sepsis
$ws.Range("B52").Value = 0.33
$ws.Range("C52").ClearContents()

# Restore the view: scroll so row 7 is at the top and select B50:B52
# (matches the sheetView/selection recorded when the author saved the
# file after making the edits above).
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B50:B52").Select()
